# Insert a new data row at row 16 (pushing existing rows 16-132 down to 17-133)
# and populate it with the new weekly price-report entry for Orégano.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Insert()

$ws.Cells.Item(16,1).Value = 6
$ws.Cells.Item(16,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16,3).Value = "Metropolitana"
$fecha = Get-Date -Year 2021 -Month 12 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(16,4).Value = $fecha
$ws.Cells.Item(16,5).Value = 13
$ws.Cells.Item(16,6).Value = 100112029
$ws.Cells.Item(16,7).Value = "Orégano"
$ws.Cells.Item(16,8).Value = "Sin especificar"
$ws.Cells.Item(16,9).Value = "Primera"
$ws.Cells.Item(16,10).Value = 35
$ws.Cells.Item(16,11).Value = 9000
$ws.Cells.Item(16,12).Value = 10000
$ws.Cells.Item(16,13).Value = 9457
$ws.Cells.Item(16,14).Value = "$/docena de atados"
$ws.Cells.Item(16,15).Value = "Región Metropolitana"
$ws.Cells.Item(16,16).Value = 3152
$ws.Cells.Item(16,17).Value = 3
$ws.Cells.Item(16,18).Value = "Hortaliza"
